# "first change in deploy mac"
# - Adds an "August 2016" label (D1/E1) to the Sheet tab
# - Updates the visit counts for the existing 4 people
# - Appends 9 new people rows (6-14) with UID / Name / visit count
# - Clears the (broken) SUM formula on Monthly_STAT!C3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# New header labels in D1 / E1 - format as Text first so Excel doesn't
# reinterpret the "August 2016" string as a date serial.
$ws.Range("D1:E1").NumberFormat = "@"
$ws.Range("D1").Value = "August 2016"
$ws.Range("E1").Value = "August 2016"

# Updated visit counts for the existing rows
$ws.Range("C2").Value = 24
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 12

# New rows of data
$newRows = @(
    @("CFD8A78940", "Nicole Latta", 3),
    @("8FD8940860", "David schachner", 3),
    @("0FD8A9BD80", "Jaspreet Kaur", 1),
    @("CFD8A9E820", "Richard Pusateri", 2),
    @("0FD8A82F40", "Mario Regino", 1),
    @("CFD8AA9A20", "Randell Holland", 1),
    @("CFD89A9C80", "Bryan Williams", 1),
    @("4FD8A85BA0", "Hojin Euam", 2),
    @("0FD8A290A0", "Paul Fabro", 1)
)

$r = 6
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Monthly_STAT: the SUM formula pointed at a stale "sheet" reference
# (#REF!) - clear it out, leaving a blank cell.
$stat = $wb.Worksheets.Item("Monthly_STAT")
$stat.Range("C3").ClearContents()

Write-Output "edits applied"
